$d = $word.ActiveDocument

$old = "Dates de la campanya Constel" + [char]0x00B7 + "laci" + [char]0x00F3 + " de Perseu 2022: 16-25 de gener, 7-16 de novembre, 6-15 de desembre"
$new = "Dates de la campanya 2022 en qu" + [char]0x00E8 + " usem la constel" + [char]0x00B7 + "laci" + [char]0x00F3 + ", Constel" + [char]0x00B7 + "laci" + [char]0x00F3 + " de Perseu 16-25 de gener, 7-16 de novembre, 6-15 de desembre"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
